$wb = $excel.ActiveWorkbook

# Rename sheets for excel/csv consistency
$wsFields = $wb.Worksheets.Item("field_mapping")
$wsFields.Name = "fields"

$wsValues = $wb.Worksheets.Item("value_mapping")
$wsValues.Name = "values"

# Move the active/selected tab from "fields" to "values"
$wsFields.Select()
$wsValues.Select()
